$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (row number => Fecha, Volumen, Precio minimo/maximo/promedio/Kg, Origen)
$rows = @{
    2  = @{ D = 44574; M = 200; N = 3000; O = 3000; P = 3000; S = 3000; R = "Región de La Araucanía" }
    3  = @{ D = 44567; M = 80;  N = 2400; O = 2400; P = 2400; S = 2400; R = "Región de La Araucanía" }
    4  = @{ D = 44176; M = 20;  N = 3000; O = 3000; P = 3000; S = 3000; R = "Región de O'Higgins" }
    5  = @{ D = 44998; M = 20;  N = 2500; O = 2500; P = 2500; S = 2500; R = "Región de La Araucanía" }
    6  = @{ D = 44999; M = 25;  N = 2500; O = 2500; P = 2500; S = 2500; R = "Región de La Araucanía" }
    8  = @{ D = 44323; M = 20;  N = 3200; O = 3200; P = 3200; S = 3200; R = "Región de La Araucanía" }
    9  = @{ D = 44215; M = 65;  N = 2800; O = 2800; P = 2800; S = 2800; R = "Región de La Araucanía" }
    10 = @{ D = 44175; M = 40;  N = 5000; O = 5000; P = 5000; S = 5000; R = "Provincia de Curicó" }
    11 = @{ D = 44551; M = 120; N = 4500; O = 4500; P = 4500; S = 4500; R = "Región de O'Higgins" }
    12 = @{ D = 44214; M = 50;  N = 1800; O = 1800; P = 1800; S = 1800; R = "Región de La Araucanía" }
    13 = @{ D = 44616; M = 200; N = 3200; O = 3200; P = 3200; S = 3200; R = "Región de La Araucanía" }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
}
